$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row tweaks -------------------------------------------------
# "Day Enrolled" -> "Date Enrolled (MM/DD/YYYY)"
$ws.Range("K1").Value = "Date Enrolled (MM/DD/YYYY)"

# Drop the trailing "No Of Modules Undertaking" header column and the
# stray hashed placeholder value that used to sit in A2.
$ws.Range("M1").Clear()
$ws.Range("A2").Clear()

# --- New sample/data-entry rows -----------------------------------------
# Email Address column gets a blue font (data-entry hint colour) for a
# few rows, Date Of Birth and Date Enrolled columns get a date number
# format so typed dates render as mm/dd/yy.
$ws.Range("D2:D4").Font.Color = 16711680
$ws.Range("G2:G5").NumberFormat = "mm/dd/yy"
$ws.Range("K2:K5").NumberFormat = "mm/dd/yy"

# --- Column width touch-ups ----------------------------------------------
# Email Address and Date Enrolled columns grow a bit to fit the new
# longer header text / blue highlighted sample cells.
$ws.Columns.Item(4).ColumnWidth = 25.25
$ws.Columns.Item(11).ColumnWidth = 25.59

# --- Selection / scroll position -----------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
$ws.Range("K4").Select() | Out-Null
